# Negate every numeric value in column E ("Block" frame data) across all
# data rows of the sheet. Blank / non-numeric cells are left untouched, and
# zero values are naturally unaffected by negation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)  # Column E
    $val = $cell.Value2
    if ($val -is [double]) {
        $cell.Value2 = -1 * $val
    }
}
